$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The crawl re-ran and a handful of products swapped positions in the sheet;
# move the affected data rows (columns A:N) into their new positions using a
# scratch row as temp storage (clearing destinations first, since copying an
# empty source cell otherwise leaves stale data behind), then clear the
# scratch row when done.

function Move-Row($srcRow, $dstRow) {
    $ws.Range("A${dstRow}:N${dstRow}").Clear()
    $ws.Range("A${srcRow}:N${srcRow}").Copy($ws.Range("A${dstRow}:N${dstRow}"))
}

# Simple pairwise swaps
$swapPairs = @(
    @(9, 10),
    @(32, 33),
    @(35, 36),
    @(42, 43),
    @(48, 49),
    @(50, 51),
    @(85, 86)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    Move-Row $r1 200
    Move-Row $r2 $r1
    Move-Row 200 $r2
    $ws.Range("A200:N200").Clear()
}

# Rotation of rows 25, 26, 27: new25 = old26, new26 = old27, new27 = old25
Move-Row 25 200
Move-Row 26 25
Move-Row 27 26
Move-Row 200 27
$ws.Range("A200:N200").Clear()

# Rotation of rows 68, 69, 70: new68 = old70, new69 = old68, new70 = old69
Move-Row 70 200
Move-Row 69 70
Move-Row 68 69
Move-Row 200 68
$ws.Range("A200:N200").Clear()

# Refresh the crawl timestamp column for every data row.
$ws.Range("O2:O88").Value = "2022-07-29 20:57:51"

Write-Host "done"
